$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$emails = @(
    "HENNOWS001@retail.spar.co.za",
    "HENNOWS002@retail.spar.co.za",
    "HENNOWS003@retail.spar.co.za",
    "HENNOWS004@retail.spar.co.za",
    "HENNOWS005@retail.spar.co.za",
    "HENNOWS006@retail.spar.co.za",
    "HENNOWS007@retail.spar.co.za",
    "HENNOWS008@retail.spar.co.za"
)

for ($i = 0; $i -lt $emails.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $emails[$i]
}
